# Scheduled runner update: refresh currentAveragePrice/Leve profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Profits" sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 2111.913
$ws.Range("I15").Value = 2111.913
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6335.739
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6166.739
# row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = ""
# row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = ""
# row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
# row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
# row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
# row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""
# row 112
$ws.Range("H112").Value = 3211.8147
$ws.Range("I112").Value = 995
$ws.Range("J112").Value = 3488.9167
$ws.Range("K112").Value = 2985
$ws.Range("L112").Value = 10466.7501
$ws.Range("M112").Value = -1877
$ws.Range("N112").Value = -12682.7501
# row 132
$ws.Range("H132").Value = 1742
$ws.Range("I132").Value = 1278
$ws.Range("J132").Value = 4990
$ws.Range("K132").Value = 3834
$ws.Range("L132").Value = 14970
$ws.Range("M132").Value = -1304
# row 137
$ws.Range("H137").Value = 2178.3572
$ws.Range("I137").Value = 2215.1538
$ws.Range("J137").Value = 1700
$ws.Range("K137").Value = 6645.4614
$ws.Range("L137").Value = 5100
$ws.Range("M137").Value = -4095.4614
$ws.Range("N137").Value = -10200
# row 138
$ws.Range("H138").Value = 2575.0557
$ws.Range("I138").Value = 1291.3
$ws.Range("J138").Value = 4179.75
$ws.Range("K138").Value = 3873.9
$ws.Range("L138").Value = 12539.25
$ws.Range("M138").Value = 1266.1
$ws.Range("N138").Value = -22819.25

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1944.95
$ws.Range("I32").Value = 1791.8379
$ws.Range("J32").Value = 3833.3333
$ws.Range("K32").Value = 1791.8379
$ws.Range("L32").Value = 3833.3333
$ws.Range("M32").Value = -1504.8379
# row 61
$ws.Range("H61").Value = 2329.65
$ws.Range("I61").Value = 1435.2142
$ws.Range("J61").Value = 4416.6665
$ws.Range("K61").Value = 1435.2142
$ws.Range("L61").Value = 4416.6665
$ws.Range("M61").Value = -1223.2142
# row 132
$ws.Range("H132").Value = 2701.077
$ws.Range("I132").Value = 1743.125
$ws.Range("J132").Value = 4233.8
$ws.Range("K132").Value = 5229.375
$ws.Range("L132").Value = 12701.4
$ws.Range("M132").Value = -2699.375
$ws.Range("N132").Value = -17761.4
# row 136
$ws.Range("H136").Value = 2329.65
$ws.Range("I136").Value = 1435.2142
$ws.Range("J136").Value = 4416.6665
$ws.Range("K136").Value = 4305.642599999999
$ws.Range("L136").Value = 13249.9995
$ws.Range("M136").Value = -1755.642599999999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 1499.75
$ws.Range("I99").Value = 1666.3334
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1666.3334
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -168.3334
$ws.Range("N99").Value = -3996
# row 102
$ws.Range("H102").Value = 1111
$ws.Range("I102").Value = 1111
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1111
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 2134

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = ""
# row 31
$ws.Range("H31").Value = 2948.8667
$ws.Range("I31").Value = 1154.5
$ws.Range("J31").Value = 4999.5713
$ws.Range("K31").Value = 1154.5
$ws.Range("L31").Value = 4999.5713
$ws.Range("M31").Value = -859.5
$ws.Range("N31").Value = -5589.5713
# row 34
$ws.Range("H34").Value = 2948.8667
$ws.Range("I34").Value = 1154.5
$ws.Range("J34").Value = 4999.5713
$ws.Range("K34").Value = 1154.5
$ws.Range("L34").Value = 4999.5713
$ws.Range("M34").Value = -952.5
$ws.Range("N34").Value = -5403.5713
# row 58
$ws.Range("H58").Value = 1012
$ws.Range("I58").Value = 1012
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1012
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -809
$ws.Range("N58").Value = ""
# row 62
$ws.Range("H62").Value = 2921.2
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 5006
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 5006
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -6254
# row 65
$ws.Range("H65").Value = 2921.2
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 5006
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 25030
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -31270
# row 97
$ws.Range("H97").Value = 16000
$ws.Range("I97").Value = 11500
$ws.Range("J97").Value = 25000
$ws.Range("K97").Value = 11500
$ws.Range("L97").Value = 25000
$ws.Range("M97").Value = -10509
$ws.Range("N97").Value = -26982
# row 104
$ws.Range("H104").Value = 102500
$ws.Range("I104").Value = 80000
$ws.Range("J104").Value = 125000
$ws.Range("K104").Value = 80000
$ws.Range("L104").Value = 125000
$ws.Range("M104").Value = -77379
$ws.Range("N104").Value = -130242
# row 132
$ws.Range("H132").Value = 3119
$ws.Range("I132").Value = 3165.5557
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 9496.667099999999
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -6966.667099999999
# row 136
$ws.Range("H136").Value = 1012
$ws.Range("I136").Value = 1012
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3036
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -486
$ws.Range("N136").Value = ""

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 55
$ws.Range("H55").Value = 11404.875
$ws.Range("I55").Value = 414.33334
$ws.Range("J55").Value = 17999.2
$ws.Range("K55").Value = 1243.00002
$ws.Range("L55").Value = 53997.60000000001
$ws.Range("M55").Value = -1066.00002
$ws.Range("N55").Value = -54351.60000000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 359.75
$ws.Range("I97").Value = 351.14285
$ws.Range("J97").Value = 420
$ws.Range("K97").Value = 351.14285
$ws.Range("L97").Value = 420
$ws.Range("M97").Value = 144.85715
# row 132
$ws.Range("H132").Value = 2015.1538
$ws.Range("I132").Value = 2058.0833
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 6174.249899999999
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -3644.249899999999

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 2526
$ws.Range("I46").Value = 1764.7693
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1764.7693
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1576.7693
$ws.Range("N46").Value = -5376
# row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
# row 136
$ws.Range("H136").Value = 959.3333
$ws.Range("I136").Value = 939
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2817
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -267
$ws.Range("N136").Value = -8100

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 70
$ws.Range("H70").Value = 20000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630
# row 73
$ws.Range("H73").Value = 20000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184
# row 132
$ws.Range("H132").Value = 2177.639
$ws.Range("I132").Value = 984.16
$ws.Range("J132").Value = 4890.091
$ws.Range("K132").Value = 2952.48
$ws.Range("L132").Value = 14670.273
$ws.Range("M132").Value = -422.48
$ws.Range("N132").Value = -19730.273
# row 136
$ws.Range("H136").Value = 983.625
$ws.Range("I136").Value = 973.0454999999999
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 2919.1365
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = -369.1364999999996
